# Update Name of Algo
# Apply corrected imputed values (columns B and C) for the KNN result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.091
$ws.Range("B3").Value = 5.305
$ws.Range("B14").Value = 5.735
$ws.Range("B16").Value = 5.394
$ws.Range("C18").Value = -11.742
$ws.Range("B21").Value = 9.270999999999999
$ws.Range("B23").Value = 8.000999999999999
$ws.Range("C24").Value = -12.704
$ws.Range("B25").Value = 6.176999999999999
$ws.Range("C25").Value = -13.959
$ws.Range("B26").Value = 6.203
$ws.Range("C27").Value = -13.203
$ws.Range("B29").Value = 5.319999999999999
$ws.Range("C30").Value = -13.165
$ws.Range("C31").Value = -12.804
$ws.Range("C39").Value = -12.702
$ws.Range("B40").Value = 9.4
$ws.Range("C42").Value = -12.752
$ws.Range("C48").Value = -11.22
$ws.Range("C51").Value = -11.02
$ws.Range("C52").Value = -11.303
$ws.Range("B53").Value = 6.293
$ws.Range("C55").Value = -13.748
$ws.Range("C56").Value = -13.222
$ws.Range("B57").Value = 5.575
$ws.Range("C57").Value = -13.385
$ws.Range("B59").Value = 4.914
$ws.Range("C60").Value = -12.8
$ws.Range("B65").Value = 5.936000000000001
$ws.Range("B69").Value = 5.319999999999999
$ws.Range("C73").Value = -12.886
$ws.Range("C74").Value = -12.453
$ws.Range("B79").Value = 5.566
$ws.Range("B83").Value = 5.523999999999999
$ws.Range("C90").Value = -13.298
$ws.Range("B91").Value = 5.545
$ws.Range("C92").Value = -11.087
$ws.Range("B93").Value = 5.409000000000001
$ws.Range("B100").Value = 5.893000000000001
